$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = "70.053.51"
$ws.Range("D2").Style = "Normal"
$ws.Range("E2").Value = "  -1.32%  "

$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = "3.576.96"
$ws.Range("D3").Style = "Normal"
$ws.Range("E3").Value = "  -1.94%  "

$ws.Range("E4").Value = "  +0.09%  "

$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "577.50"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = "  -2.84%  "

$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "186.14"
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = "  -5.05%  "

$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "3.573.76"
$ws.Range("D7").Style = "Normal"
$ws.Range("E7").Value = "  -1.90%  "

$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "0.618"
$ws.Range("D8").Style = "Normal"
$ws.Range("E8").Value = "  -4.74%  "

$ws.Range("E9").Value = "  +0.07%  "

$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "0.183"
$ws.Range("D10").Style = "Normal"
$ws.Range("E10").Value = "  -0.84%  "

$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "0.650"
$ws.Range("D11").Style = "Normal"
$ws.Range("E11").Value = "  -4.58%  "

$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "54.92"
$ws.Range("D12").Style = "Normal"
$ws.Range("E12").Value = "  -5.91%  "

$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "0.0000304"
$ws.Range("D13").Style = "Normal"
$ws.Range("E13").Value = "  +2.56%  "

$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "9.53"
$ws.Range("D14").Style = "Normal"
$ws.Range("E14").Value = "  -4.74%  "

$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "4.152.77"
$ws.Range("D15").Style = "Normal"
$ws.Range("E15").Value = "  -1.86%  "

$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "19.65"
$ws.Range("D16").Style = "Normal"
$ws.Range("E16").Value = "  -3.75%  "

$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "3.583.97"
$ws.Range("D17").Style = "Normal"
$ws.Range("E17").Value = "  -1.72%  "

$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "70.008.51"
$ws.Range("D18").Style = "Normal"
$ws.Range("E18").Value = "  -1.31%  "

$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "12.56"
$ws.Range("D19").Style = "Normal"
$ws.Range("E19").Value = "  -1.75%  "

$ws.Range("E20").Value = "  -0.93%  "

$ws.Range("E21").Value = "  -3.27%  "

$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "490.69"
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = "  +0.15%  "

$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "19.08"
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = "  -1.84%  "

$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "4.95"
$ws.Range("D24").Style = "Normal"
$ws.Range("E24").Value = "  -5.84%  "

$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "4.39"
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = "  -2.44%  "

$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "94.89"
$ws.Range("D26").Style = "Normal"
$ws.Range("E26").Value = "  +3.62%  "

$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "11.75"
$ws.Range("D27").Style = "Normal"
$ws.Range("E27").Value = "  +2.52%  "

$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "2.96"
$ws.Range("D28").Style = "Normal"
$ws.Range("E28").Value = "  -6.83%  "

$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "9.32"
$ws.Range("D29").Style = "Normal"
$ws.Range("E29").Value = "  -3.30%  "

$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "7.77"
$ws.Range("D30").Style = "Normal"
$ws.Range("E30").Value = "  -2.15%  "

$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "31.57"
$ws.Range("D31").Style = "Normal"
$ws.Range("E31").Value = "  -3.95%  "

$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "12.08"
$ws.Range("D32").Style = "Normal"
$ws.Range("E32").Value = "  -1.84%  "

$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "65.84"
$ws.Range("D33").Style = "Normal"
$ws.Range("E33").Value = "  -1.11%  "

$ws.Range("E34").Value = "  -6.59%  "

$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "579.95"
$ws.Range("D35").Style = "Normal"
$ws.Range("E35").Value = "  -6.11%  "

$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "3.26"
$ws.Range("D36").Style = "Normal"
$ws.Range("E36").Value = "  +14.26%  "

$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "0.417"
$ws.Range("D37").Style = "Normal"
$ws.Range("E37").Value = "  +0.82%  "

$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "38.78"
$ws.Range("D38").Style = "Normal"
$ws.Range("E38").Value = "  -3.94%  "

$ws.Range("E39").Value = "  +0.09%  "

$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "0.0₃0790"
$ws.Range("D40").Style = "Normal"
$ws.Range("E40").Value = "  -5.32%  "

$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "3.47"
$ws.Range("D41").Style = "Normal"
$ws.Range("E41").Value = "  -3.40%  "

$ws.Range("E42").Value = "  -2.27%  "

$ws.Range("E43").Value = "  -9.41%  "

$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "3.07"
$ws.Range("D44").Style = "Normal"
$ws.Range("E44").Value = "  -3.34%  "

$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "3.206.09"
$ws.Range("D45").Style = "Normal"
$ws.Range("E45").Value = "  -3.83%  "

$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "0.0443"
$ws.Range("D46").Style = "Normal"
$ws.Range("E46").Value = "  -3.52%  "

$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "3.47"
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = "  +4.73%  "

$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "1.64"
$ws.Range("D48").Style = "Normal"
$ws.Range("E48").Value = "  +35.55%  "

$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "9.58"
$ws.Range("D49").Style = "Normal"
$ws.Range("E49").Value = "  -0.79%  "

$ws.Range("E50").Value = "  -2.87%  "

$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "1.00"
$ws.Range("D51").Style = "Normal"
$ws.Range("E51").Value = "  +0.05%  "
